# Update odds values on "Sheet1" for rows 3 and 6 of the
# Jogos_da_Semana_FlashScore_2024-11-13 workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 updates
$ws.Range("Q3").Value = 1.83
$ws.Range("R3").Value = 2.03

# Row 6 updates
$ws.Range("G6").Value = 1.4
$ws.Range("H6").Value = 3.75
$ws.Range("J6").Value = 1.87
$ws.Range("L6").Value = 7.9
$ws.Range("P6").Value = 2.85
$ws.Range("Q6").Value = 1.93
$ws.Range("R6").Value = 1.7
$ws.Range("W6").Value = 5.3
$ws.Range("AD6").Value = 7.8
$ws.Range("AE6").Value = 23
$ws.Range("AF6").Value = 150
$ws.Range("AH6").Value = 18.5
$ws.Range("AL6").Value = 150
$ws.Range("AM6").Value = 120
$ws.Range("AO6").Value = 6.3
$ws.Range("AQ6").Value = 18.5
$ws.Range("AU6").Value = 8.25
$ws.Range("AV6").Value = 90
$ws.Range("AX6").Value = 60
$ws.Range("AY6").Value = 55
$ws.Range("AZ6").Value = 500
$ws.Range("BA6").Value = 450
